$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet1 - "Test Series System" (columns A:L, formula uses B/G -> L)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Test Series System")

# Row 3 updates
$ws1.Range("B3").Value = 1
$ws1.Range("F3").Value = 1
$ws1.Range("G3").Value = 1
$ws1.Range("K3").Value = 1

# Row 4 updates
$ws1.Range("E4").Value = 3
$ws1.Range("F4").Value = 1
$ws1.Range("G4").Value = 1
$ws1.Range("K4").Value = 1

# Row 5 updates
$ws1.Range("B5").Value = 1
$ws1.Range("C5").Value = 3
$ws1.Range("E5").Value = 3
$ws1.Range("F5").Value = 1
$ws1.Range("G5").Value = 1
$ws1.Range("K5").Value = 1

# New row 6 - copy formatting from row 5 first, then set values/formula
$ws1.Range("A5:L5").Copy()
$ws1.Range("A6:L6").PasteSpecial(-4122)
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 1
$ws1.Range("H6").Value = 3
$ws1.Range("I6").Value = 3
$ws1.Range("J6").Value = 3
$ws1.Range("K6").Value = 1
$ws1.Range("L6").Formula = "=IF(B6 = G6, 1, 0)"

# Extend conditional formatting range
$ws1.Range("L2:L5").FormatConditions.Item(1).ModifyAppliesToRange($ws1.Range("L2:L6"))

# ------------------------------------------------------------------
# Sheet2 - "Sensed high_rel_comp1 History" (columns A:K)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sensed high_rel_comp1 History")

# Row 3 updates
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0

# Row 4 updates
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 0

# Row 5 updates
$ws2.Range("B5").Value = 3
$ws2.Range("C5").Value = 1

# New row 6
$ws2.Range("A5:K5").Copy()
$ws2.Range("A6:K6").PasteSpecial(-4122)
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 1
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 3
$ws2.Range("G6").Value = 3
$ws2.Range("H6").Value = 3
$ws2.Range("I6").Value = 3
$ws2.Range("J6").Formula = "=IF(B6 = F6, 1, 0)"
$ws2.Range("K6").Formula = "=MODE(C6:E6)"

# Extend conditional formatting ranges
$ws2.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("J2:J6"))
$ws2.Range("K2:K5").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("K2:K6"))

# ------------------------------------------------------------------
# Sheet3 - "Sensed high_rel_comp2 History" (columns A:K)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sensed high_rel_comp2 History")

# Row 3 updates
$ws3.Range("C3").Value = 1
$ws3.Range("E3").Value = 0

# Row 4 updates
$ws3.Range("C4").Value = 1
$ws3.Range("E4").Value = 0

# Row 5 updates
$ws3.Range("C5").Value = 1

# New row 6
$ws3.Range("A5:K5").Copy()
$ws3.Range("A6:K6").PasteSpecial(-4122)
$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = 0
$ws3.Range("C6").Value = 0
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 0
$ws3.Range("F6").Value = 3
$ws3.Range("G6").Value = 3
$ws3.Range("H6").Value = 3
$ws3.Range("I6").Value = 3
$ws3.Range("J6").Formula = "=IF(B6 = F6, 1, 0)"
$ws3.Range("K6").Formula = "=MODE(C6:E6)"

# Extend conditional formatting ranges
$ws3.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("J2:J6"))
$ws3.Range("K2:K5").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("K2:K6"))

# ------------------------------------------------------------------
# Sheet4 - "Sensed low_rel_comp1 History" (columns A:K)
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sensed low_rel_comp1 History")

# Row 3 updates
$ws4.Range("E3").Value = 0

# Row 4 updates
$ws4.Range("B4").Value = 3
$ws4.Range("D4").Value = 0
$ws4.Range("H4").Value = 3

# Row 5 updates
$ws4.Range("B5").Value = 3
$ws4.Range("H5").Value = 3

# New row 6
$ws4.Range("A5:K5").Copy()
$ws4.Range("A6:K6").PasteSpecial(-4122)
$ws4.Range("A6").Value = 4
$ws4.Range("B6").Value = 0
$ws4.Range("C6").Value = 0
$ws4.Range("D6").Value = 0
$ws4.Range("E6").Value = 0
$ws4.Range("F6").Value = 3
$ws4.Range("G6").Value = 3
$ws4.Range("H6").Value = 3
$ws4.Range("I6").Value = 3
$ws4.Range("J6").Formula = "=IF(B6 = F6, 1, 0)"
$ws4.Range("K6").Formula = "=MODE(C6:E6)"

# Extend conditional formatting ranges
$ws4.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($ws4.Range("J2:J6"))
$ws4.Range("K2:K5").FormatConditions.Item(1).ModifyAppliesToRange($ws4.Range("K2:K6"))

# ------------------------------------------------------------------
# Sheet5 - "Sensed low_rel_comp2 History" (columns A:K)
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sensed low_rel_comp2 History")

# Row 3 updates
$ws5.Range("B3").Value = 1
$ws5.Range("C3").Value = 1
$ws5.Range("F3").Value = 1
$ws5.Range("G3").Value = 1
$ws5.Range("H3").Value = 1
$ws5.Range("I3").Value = 1

# Row 4 updates
$ws5.Range("B4").Value = 1
$ws5.Range("F4").Value = 1
$ws5.Range("G4").Value = 1
$ws5.Range("H4").Value = 1
$ws5.Range("I4").Value = 1

# Row 5 updates
$ws5.Range("B5").Value = 1
$ws5.Range("E5").Value = 0
$ws5.Range("F5").Value = 1
$ws5.Range("G5").Value = 1
$ws5.Range("H5").Value = 1
$ws5.Range("I5").Value = 1

# New row 6
$ws5.Range("A5:K5").Copy()
$ws5.Range("A6:K6").PasteSpecial(-4122)
$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = 0
$ws5.Range("C6").Value = 0
$ws5.Range("D6").Value = 1
$ws5.Range("E6").Value = 0
$ws5.Range("F6").Value = 1
$ws5.Range("G6").Value = 1
$ws5.Range("H6").Value = 0
$ws5.Range("I6").Value = 1
$ws5.Range("J6").Formula = "=IF(B6 = F6, 1, 0)"
$ws5.Range("K6").Formula = "=MODE(C6:E6)"

# Extend conditional formatting ranges
$ws5.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($ws5.Range("J2:J6"))
$ws5.Range("K2:K5").FormatConditions.Item(1).ModifyAppliesToRange($ws5.Range("K2:K6"))
